$wb = $excel.ActiveWorkbook

# --- Sheet "Productdata": update AverageDemand (column G) for products 1-4 (rows 2-5) ---
$wsProduct = $wb.Worksheets.Item("Productdata")
$wsProduct.Range("G2").Value = 49
$wsProduct.Range("G3").Value = 21
$wsProduct.Range("G4").Value = 35
$wsProduct.Range("G5").Value = 70

# Column H (StandardDevDemands) holds empty placeholder cells (blank shared
# strings) in the source file. Re-assert them as blank so the round trip
# through the COM layer doesn't coerce them into a stray text value.
for ($r = 2; $r -le 11; $r++) {
    $wsProduct.Cells.Item($r, 8).Value = ""
}

# --- Sheet "ForecastedAverageDemand": fill in demand for the last 3 periods (rows 9-11) ---
$wsDemand = $wb.Worksheets.Item("ForecastedAverageDemand")

$wsDemand.Range("B9").Value = 70
$wsDemand.Range("C9").Value = 30
$wsDemand.Range("D9").Value = 50
$wsDemand.Range("E9").Value = 100

$wsDemand.Range("B10").Value = 70
$wsDemand.Range("C10").Value = 30
$wsDemand.Range("D10").Value = 50
$wsDemand.Range("E10").Value = 100

$wsDemand.Range("B11").Value = 70
$wsDemand.Range("C11").Value = 30
$wsDemand.Range("D11").Value = 50
$wsDemand.Range("E11").Value = 100

# --- Sheet "ForcastedStandardDeviation": fill in matching standard deviations for rows 9-11 ---
$wsStdDev = $wb.Worksheets.Item("ForcastedStandardDeviation")

$wsStdDev.Range("B9").Value = 7.166424999999998
$wsStdDev.Range("C9").Value = 3.071324999999999
$wsStdDev.Range("D9").Value = 5.118874999999999
$wsStdDev.Range("E9").Value = 10.23775

$wsStdDev.Range("B10").Value = 8.1997825
$wsStdDev.Range("C10").Value = 3.5141925
$wsStdDev.Range("D10").Value = 5.856987499999999
$wsStdDev.Range("E10").Value = 11.713975

$wsStdDev.Range("B11").Value = 9.129804249999998
$wsStdDev.Range("C11").Value = 3.912773249999999
$wsStdDev.Range("D11").Value = 6.521288749999998
$wsStdDev.Range("E11").Value = 13.0425775
